# Update quizvragen via Admin
#
# Sheet "DC": record the image that was generated for the existing last
# question (row 6, column L = image_url) and append a brand-new "mc"
# question in row 7, which also already has its image_url filled in.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DC")

# Row 6 gained its image_url value.
$ws.Range("L6").Value = "https://raw.githubusercontent.com/onomatorHanze/didactic-octo-spork/main/data/images/DC_new_1763126867.png"

# New row 7: a fresh "mc" question. Most metadata columns are still blank
# for this question (id, topic, explanation, image_path, formula_latex,
# tags, difficulty), same as the blank placeholders used for rows 5 and 6.
$ws.Range("A7").Style = "Normal"
$ws.Range("B7").Value = "mc"
$ws.Range("C7").Style = "Normal"
$ws.Range("D7").Value = "Bereken hoeveel regen er valt. "
$ws.Range("E7").Value = "['']"
$ws.Range("F7").Value = 0
$ws.Range("G7").Style = "Normal"
$ws.Range("H7").Style = "Normal"
$ws.Range("I7").Style = "Normal"
$ws.Range("J7").Style = "Normal"
$ws.Range("K7").Style = "Normal"
$ws.Range("L7").Value = "https://raw.githubusercontent.com/onomatorHanze/didactic-octo-spork/main/data/images/DC_new_1763130529.png"
